# Contest 8 SRH vs MI
# Fill in points data for Match 8 (row 20, SRH vs MI) and add the next
# set of match fixtures (Matches 9-26, rows 21-38) with their team names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Match 8 (row 20, "SRH vs MI") player scores ---
$ws.Range("E20").Value = 50
$ws.Range("H20").Value = 30
$ws.Range("K20").Value = 40
$ws.Range("N20").Value = 60
$ws.Range("Q20").Value = 70
$ws.Range("T20").Value = 0
$ws.Range("W20").Value = 100
$ws.Range("Z20").Value = 20

# --- Upcoming match fixtures for rows 21-38 (Matches 9-26) ---
$upcomingMatches = @(
    "RR vs DC",
    "RCB vs KKR",
    "LSG vs PBKS",
    "GT vs SRH",
    "DC vs CSK",
    "MI vs RR",
    "RCB vs LSG",
    "DC vs KKR",
    "GT vs PBKS",
    "SRH vs CSK",
    "RR vs RCB",
    "MI vs DC",
    "LSG vs GT",
    "CSK vs KKR",
    "PBKS vs SRH",
    "RR vs GT",
    "MI vs RCB",
    "LSG vs DC"
)

$row = 21
foreach ($matchName in $upcomingMatches) {
    $ws.Range("C$row").Value = $matchName
    $row++
}

$excel.Calculate()
